$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.451.89'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').Value = '1.629.39'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '304.83'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3765'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3652'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '51.73'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08215'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.42%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.223'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.73%  '
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.41'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.550'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001250'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.243'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.12%  '
$ws.Range('D17').Value = '1.631.93'
$ws.Range('E17').Value = '  -0.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '94.04'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06982'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.73'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.453'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.63%  '
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.72'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.76%  '
$ws.Range('D24').Value = '23.450.76'
$ws.Range('E24').Value = '  -0.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.171'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.462'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.38'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '150.45'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.312'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.68'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.37%  '
$ws.Range('D31').Value = '1.811.96'
$ws.Range('E31').Value = '  -0.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.263'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.795'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.017'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '10.78'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02785'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2526'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.87%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.08768'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.07135'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.034'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.7043'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.349'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.36'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.68%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.27'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.91%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6550'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.325'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.001'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.986'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08017'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.203'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '125.45'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.17%  '
